$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Fullstack Developer"
$ws.Range("C4").Value = "jghvfh"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 17
